# "Generate Report for handback" — mark the two content files as handed
# back (in sync with en-US) on both the zh-cn and de-de report sheets:
#   - Status -> "Handed back: in sync with en-US"
#   - Latest Target File  (E) <- Source File Name (A)   [new hyperlink]
#   - Latest Handback File (F) <- Latest Handoff File (C) [new hyperlink]
#   - Latest Handback DateTime (G) <- new timestamp

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Overview sheet ---------------------------------------------------------
# (shares the "Ready for handoff" string with zh-cn/de-de; once that text
# is no longer used anywhere it is updated in place, so these cells move
# to the new status text along with the per-language report rows below)
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# zh-cn sheet -----------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $newStatus
$ws.Range("E2").Value = $ws.Range("A2").Value
$ws.Range("F2").Value = $ws.Range("C2").Value
$ws.Range("G2").Value = "2016-01-26 12:29:19"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5ca88405635ef0c3bb21fdd52673af46fcb9cee7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.zh-cn.xlf")

$ws.Range("B3").Value = $newStatus
$ws.Range("E3").Value = $ws.Range("A3").Value
$ws.Range("F3").Value = $ws.Range("C3").Value
$ws.Range("G3").Value = "2016-01-26 12:29:19"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5ca88405635ef0c3bb21fdd52673af46fcb9cee7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.zh-cn.xlf")

# de-de sheet -------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = $newStatus
$ws.Range("E2").Value = $ws.Range("A2").Value
$ws.Range("F2").Value = $ws.Range("C2").Value
$ws.Range("G2").Value = "2016-01-26 12:29:39"
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9f5d0492f4d5d912b973678a82e27543a2d191b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf", "", "", "26508feb-40e2-4f19-bcc0-5b63fdbf85d6.7b8a59b1bd7f30c0199197e8a72a59e3b09d092e.de-de.xlf")

$ws.Range("B3").Value = $newStatus
$ws.Range("E3").Value = $ws.Range("A3").Value
$ws.Range("F3").Value = $ws.Range("C3").Value
$ws.Range("G3").Value = "2016-01-26 12:29:39"
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/8b0ced3ef4d263211639b58a4a4bccb631754d00/e2e/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9f5d0492f4d5d912b973678a82e27543a2d191b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf", "", "", "78d9e4aa-5767-4f22-ad97-fff7c64ce1d4.c347d153ca8010479a32b5df4707d7705be9c1a0.de-de.xlf")
